# Automatic update of files.
#
# 1) Bump the "Förändrad" (changed) date in column C for every data row
#    (rows 2-39) from 2023-09-15 (45184) to 2023-09-17 (45186).
# 2) Add the designation (column A value) as a friendly-name second
#    argument to the HYPERLINK() formulas that exist in row 2 and row 3
#    (columns S, T, V, W, X, Y).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 39
$newChangedDate = 45186

# --- 1) Update the "Förändrad" date column (C) for all data rows ---
for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $ws.Cells.Item($row, 3).Value = $newChangedDate
}

# --- 2) Rebuild the HYPERLINK formulas that include a friendly name ---
$hyperlinkColumns = @(
    @{ Col = "S"; Path = "artfynd"; Ext = "xlsx" },
    @{ Col = "T"; Path = "kartor"; Ext = "png" },
    @{ Col = "V"; Path = "klagomål"; Ext = "docx" },
    @{ Col = "W"; Path = "klagomålsmail"; Ext = "docx" },
    @{ Col = "X"; Path = "tillsyn"; Ext = "docx" },
    @{ Col = "Y"; Path = "tillsynsmail"; Ext = "docx" }
)

$rowsWithHyperlinks = @(2, 3)

foreach ($row in $rowsWithHyperlinks) {
    $designation = $ws.Cells.Item($row, 1).Value2

    foreach ($entry in $hyperlinkColumns) {
        $cellRef = "$($entry.Col)$row"
        $url = "https://klasma.github.io/Logging_SOTENAS/$($entry.Path)/$designation.$($entry.Ext)"
        $formula = '=HYPERLINK("' + $url + '", "' + $designation + '")'
        $ws.Range($cellRef).Formula = $formula
    }
}
